# Daily update: append a new data row (2020-04-18) to the "相談件数" sheet,
# pushing the trailing footnote row down by one, and extend the print area
# to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 84 (the existing row 84 - the footnote row - and
# everything below it shifts down to row 85). Excel copies the formatting
# of the row above (row 83) onto the freshly inserted row.
$ws.Rows.Item(84).Insert()

# Fill in the new day's figures (date 2020-04-18 = serial 43939).
$ws.Cells.Item(84, 1).Value = 43939
$ws.Cells.Item(84, 2).Value = 522
$ws.Cells.Item(84, 3).Value = 25457
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 5665

# Keep the print area in sync with the sheet's new extent.
$ws.PageSetup.PrintArea = '$A$1:$E$89'

# Re-select the bottom-right pane's active cell, matching where the new
# last row of data now sits.
$ws.Range("C85").Select()
